$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.045.83'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +0.09%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.830.89'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -0.03%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9982'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.06%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '243.31'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +0.45%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.6281'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +0.56%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.000'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +0.04%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.07489'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -1.00%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.2921'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +0.11%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '23.17'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +2.80%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07711'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -0.09%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.827.79'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +0.35%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.989'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +0.76%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.6679'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +0.64%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '82.58'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -0.09%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.000009313'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -7.44%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '5.996'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -0.58%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '29.076.56'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +0.28%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '2.084.79'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +0.24%  '
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +2.05%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '223.41'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -1.23%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.001'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +0.21%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.119'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -0.59%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.9998'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -0.01%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '159.22'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +0.71%  '
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +1.87%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.494'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +0.35%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '17.92'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +0.09%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.496'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +0.59%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.05709'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +9.85%  '
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +1.59%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.060'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +1.12%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.203'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +0.98%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.7469'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +1.33%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.844'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +0.05%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.136'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -0.28%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.669'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -1.03%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.757'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +0.04%  '
$ws.Range("B39").Value = 'VeChain'
$ws.Range("C39").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01783'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +0.02%  '
$ws.Range("B40").Value = 'Maker'
$ws.Range("C40").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.217.63'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -2.07%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '6.543'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +3.42%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.8919'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -0.36%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.000'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +0.06%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '102.12'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +0.79%  '
$ws.Range("B45").Value = 'RocketPoolETH'
$ws.Range("C45").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.981.60'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +0.13%  '
$ws.Range("B46").Value = 'BabyDogeCoin'
$ws.Range("C46").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.00000000124'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +0.28%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '65.73'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +2.69%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.07801'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +12.39%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.5080'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -0.52%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.4069'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +1.01%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '9.048'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +1.95%  '
